$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 64
$ws.Range("E9").Value = 36.11
$ws.Range("F9").Value = 37

$ws.Range("C10").Value = 128
$ws.Range("E10").Value = 19.260000000000002
$ws.Range("F10").Value = 20

$ws.Range("C11").Value = 256
$ws.Range("E11").Value = 10.56
$ws.Range("F11").Value = 11

$ws.Range("F11").Select()
